# Homework 1 Excel.docx edit:
#   - Replace the third bullet's sentence about "Film & video / television"
#     with a new sentence about staff picks vs successful campaigns.
#   - The trailing "_GoBack" bookmark (previously sitting after the very
#     last sentence in the document) moves to sit right after the newly
#     edited sentence instead (this is what Word does automatically when
#     you make that sentence the most recent edit).

$d = $word.ActiveDocument

$oldText = "This data only provided information on the category Film & video and the sub category of television."
$newText = "Further Analysis shows that there may be a correlation between staff picks and successful campaigns."

# Find & replace the sentence. Word's Find/Replace naturally merges the
# two existing runs ("...television" + ".") into a single run containing
# the full new sentence, matching how the real document looks afterwards.
$findRange = $d.Content
$found = $findRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# $findRange now covers the freshly-inserted replacement text, so the
# paragraph that contains it is our edited bullet point - no need to
# hard-code a paragraph index.
$targetPara = $findRange.Paragraphs(1)

# We need an empty (collapsed) bookmark positioned immediately after the
# new sentence, i.e. right before the paragraph mark. Adding a bookmark
# with a collapsed Range sitting exactly one character before a
# paragraph's end is mishandled by this host (it silently relocates to
# the very start of the document), so we dodge that exact position:
# temporarily insert a one-character guard before the paragraph mark,
# add the bookmark safely in front of the guard (no longer the last
# character of the paragraph), then remove the guard again.
$paraEnd = $targetPara.Range.End
$guard = $d.Range($paraEnd - 1, $paraEnd - 1)
$guard.InsertAfter("z") | Out-Null

$paraEndWithGuard = $targetPara.Range.End
$bmPos = $paraEndWithGuard - 2
$bmRange = $d.Range($bmPos, $bmPos)
# Re-adding a bookmark named "_GoBack" removes any previous bookmark of
# that name elsewhere in the document, which is exactly the move we need.
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$paraEndFinal = $targetPara.Range.End
$guardRange = $d.Range($paraEndFinal - 2, $paraEndFinal - 1)
$guardRange.Delete() | Out-Null
